$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.814.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.829.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6898"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07646"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3038"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07785"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "93.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.826.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.076"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6773"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.439"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008217"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.50%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.796.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.072.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.432"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1487"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.724"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.536"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.217"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.153"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.188"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05114"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7710"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.846"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.135"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.691"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01852"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.259.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.696"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9575"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.988"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.96%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.639"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5151"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.972.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.13%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000120"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.743"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.936"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.98%  "
